$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About": update the currency-year conversion factors / notes
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Existing 2022->2012 conversion factor becomes the 2023->2012 factor
$about.Range("A12").Value = 0.75350342301658668
$about.Range("B12").Value = "2023 to 2012 USD"

# New row: 2024->2012 conversion factor
$about.Range("A13").Value = 0.73
$about.Range("B13").Value = "2024 to 2012 USD"

# New row: annual inflation assumption after 2024, formatted as a percent
$about.Range("A14").Value = 0.03
$about.Range("A14").Style = "Percent"
$about.Range("B14").Value = "annual inflation assumption after 2024"

# Note under "Notes:" now clarifies the credit value is not inflation adjusted
$about.Range("A10").Value = "batteries. This value is not inflation adjusted."

# ---------------------------------------------------------------------------
# Sheet "BSfVBP": recompute the per-year credit values using the new
# currency-year conversion factors, with an assumed annual inflation decay
# applied from 2026 onward
# ---------------------------------------------------------------------------
$trans = $wb.Worksheets.Item("BSfVBP")

$trans.Range("D2").Formula = "=45*About!`$A`$12"
$trans.Range("E2").Formula = "=45*About!`$A`$13"
$trans.Range("F2").Formula = "=E2*(1-About!`$A`$14)"
$trans.Range("G2").Formula = "=F2*(1-About!`$A`$14)"
$trans.Range("H2").Formula = "=G2*(1-About!`$A`$14)"
$trans.Range("I2").Formula = "=H2*(1-About!`$A`$14)"
$trans.Range("J2").Formula = "=I2*(1-About!`$A`$14)"
$trans.Range("K2").Formula = "=J2*(1-About!`$A`$14)"
$trans.Range("L2").Formula = "=K2*(1-About!`$A`$14)"
$trans.Range("M2").Formula = "=L2*(1-About!`$A`$14)"

# Record the selection left on the BSfVBP sheet (it is not the active tab)
$trans.Range("F2:M2").Select()

# Switch back to the About sheet last so it remains the active/selected tab,
# with its own selection recorded
$about.Range("A11").Select()

$wb.Save()
